# PolicyBazar_TravelInsuranceData.xlsx - "Made changes in TravelInsurance excel file"
#
# The Switzerland row (row 5) travel dates were corrected:
#   Start Date: 2025/November/20 -> 2025/September/20
#   End Date:   2026/January/4   -> 2025/October/4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "2025/September/20"
$ws.Range("C5").Value = "2025/October/4"

# Leave the selection where the author left it when saving.
$ws.Range("C5").Select()
